# Apply "Added backup codes for crdc" change:
# - A2's used-up backup code is replaced by the next unused code that used
#   to sit at A11.
# - The now-consumed codes in A3, A4 and A11 are cleared.
# - A12's code is left untouched.
# - The active selection moves to A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = $ws.Range("A11").Value2

$ws.Range("A3").ClearContents()
$ws.Range("A4").ClearContents()
$ws.Range("A11").ClearContents()

# Update the selected cell to match the saved view state.
$ws.Range("A3").Select()
